$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Brittany"
$ws.Range("B5").Value = "Miller"
$ws.Range("C5").Value = "1406 SE Stark St."
$ws.Range("D5").Value = "Portland"
$ws.Range("E5").Value = "OR"
$ws.Range("F5").Value = 97214

$null = $ws.Range("C6").Select()
